$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (B/C) and percentage (E) columns: plain text assignment is safe,
# Excel will not mis-parse these as numbers.
$ws.Range('E2').Value = '  -1.31%  '
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -5.31%  '
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -4.06%  '
$ws.Range('E9').Value = '  -6.98%  '
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('E14').Value = '  -4.42%  '
$ws.Range('E15').Value = '  -1.98%  '
$ws.Range('E16').Value = '  +2.84%  '
$ws.Range('E17').Value = '  -1.50%  '
$ws.Range('E18').Value = '  -2.86%  '
$ws.Range('E19').Value = '  -2.53%  '
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('E22').Value = '  +7.18%  '
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -4.38%  '
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('E29').Value = '  -3.27%  '
$ws.Range('E30').Value = '  +3.67%  '
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('E32').Value = '  +3.32%  '
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  -3.87%  '
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('E41').Value = '  +10.54%  '
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E45').Value = '  +0.28%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('E48').Value = '  -2.64%  '
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E51').Value = '  +10.62%  '

# Price (D) column: force text format first so values like '243.75' or
# '1.0000' are stored as literal text (matching the original inlineStr
# cells) instead of being auto-converted to numbers by Excel.
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D19', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.807.73'
$ws.Range('D3').Value = '1.891.59'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D5').Value = '0.7781'
$ws.Range('D6').Value = '243.75'
$ws.Range('D8').Value = '0.3127'
$ws.Range('D9').Value = '25.27'
$ws.Range('D10').Value = '0.07168'
$ws.Range('D11').Value = '0.08053'
$ws.Range('D12').Value = '0.7635'
$ws.Range('D13').Value = '5.454'
$ws.Range('D14').Value = '1.856.91'
$ws.Range('D15').Value = '92.24'
$ws.Range('D16').Value = '6.155'
$ws.Range('D17').Value = '29.775.63'
$ws.Range('D19').Value = '243.37'
$ws.Range('D21').Value = '0.9995'
$ws.Range('D22').Value = '8.097'
$ws.Range('D23').Value = '2.115.27'
$ws.Range('D25').Value = '0.1618'
$ws.Range('D26').Value = '9.395'
$ws.Range('D27').Value = '161.61'
$ws.Range('D28').Value = '18.71'
$ws.Range('D29').Value = '2.047'
$ws.Range('D30').Value = '1.418'
$ws.Range('D31').Value = '1.548'
$ws.Range('D32').Value = '4.473'
$ws.Range('D33').Value = '4.100'
$ws.Range('D34').Value = '0.05535'
$ws.Range('D35').Value = '1.264'
$ws.Range('D37').Value = '0.9958'
$ws.Range('D38').Value = '2.617'
$ws.Range('D40').Value = '2.787'
$ws.Range('D41').Value = '1.138.55'
$ws.Range('D42').Value = '73.75'
$ws.Range('D43').Value = '0.4419'
$ws.Range('D44').Value = '5.852'
$ws.Range('D45').Value = '0.8527'
$ws.Range('D46').Value = '0.9998'
$ws.Range('D47').Value = '103.66'
$ws.Range('D48').Value = '1.884'
$ws.Range('D49').Value = '9.887'
$ws.Range('D50').Value = '7.442'
$ws.Range('D51').Value = '3.015'

# Reset style back to Normal so no extra cell style/format is introduced
# (keeps the cells stored as plain text without a lingering number format).
foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}
